# Sample Project / Main.xlsx - row 11 ("R40" rule) B11 now holds the text "1"
# instead of the label "R40" (value retyped as plain text, not a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
